$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Copy the formatting of the row above (row 30) into the new row 32 for
#    columns A:F - this reproduces the exact same cell styles (bordered,
#    12pt font for A:C, bold centered for F) used by the other data rows.
# ---------------------------------------------------------------------------
$ws.Range("A30:F30").Copy() | Out-Null
$ws.Range("A32:F32").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Column G/H use the plain style already used elsewhere on the sheet (e.g. G31)
$ws.Range("G31").Copy() | Out-Null
$ws.Range("G32").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$ws.Range("H32").PasteSpecial(-4122) | Out-Null        # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Fill in the new row's values
# ---------------------------------------------------------------------------
$ws.Range("A32").Value2 = 26
$ws.Range("B32").Value2 = "User Registration"
$ws.Range("C32").Value2 = "POST"
$ws.Range("D32").Value2 = "api/user-registration"
$ws.Range("E32").Value2 = "USP_AUTHENTICATE_MANAGEMENT"
$ws.Range("F32").Value2 = 4
$ws.Range("G32").Value2 = "User Registration"
$ws.Range("H32").Value2 = '{"TITLE":"12","NAME":"Sujata Kumari","GENDER":"FEMALE","DOB":"27-12-1997","EMAIL":"sujata@tts.com","CONTACT":"9988665533","PASSWORD":"pass","COUNTRY":"104","CITY":"67","COUNTRY_CODE":"+91","CHAPTER":"New Chapter Test","EDUCATION":"","ACTIVITY":"hi"}'

# ---------------------------------------------------------------------------
# 3. Row height matches the other data rows (15.75pt)
# ---------------------------------------------------------------------------
$ws.Rows.Item(32).RowHeight = 15.75

# ---------------------------------------------------------------------------
# 4. Update the view: scroll so row 19 is at the top and select H32, matching
#    the sheetView/selection recorded for the edited sheet.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H32").Select() | Out-Null
